$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the rate text in cell A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $ws1.Range("A1")
$oldText = $cellA1.Value2
$newText = $oldText.Replace(
    "1000 Bs = 15.24 = 64558.08 pesos",
    "1000 Bs = 15.13 = 63993.95 pesos"
)
$newText = $newText.Replace(
    "64558.08 pesos = 15.13 = 966.56 Bs",
    "63993.95 pesos = 15.09 = 978.32 Bs"
)
$cellA1.Value2 = $newText

# --- Sheet "tasas": update the rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value2 = 66.09999999999999
$ws2.Range("O10").Value2 = 4230
$ws2.Range("N12").Value2 = 4240
$ws2.Range("O12").Value2 = 64.81999999999999
